$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update query text in B2 (CasesTab row / "Survival (days)" query):
#     append an ORDER BY / LIMIT clause.
$b2 = $ws.Range("B2").Text
$ws.Range("B2").Value = $b2 + "`n order By ss.study_subject_id ASC LIMIT 100"

# --- Update query text in B3 (SamplesTab row / "Sample Procurement Method" query):
#     append an ORDER BY / LIMIT clause.
$b3 = $ws.Range("B3").Text
$ws.Range("B3").Value = $b3 + "`n order By samp.sample_id ASC LIMIT 100"

# --- Update query text in B4 (FilesTab row / file-name query):
#     replace the trailing "    order by f.file_name" line with the
#     new "     order By f.file_name ASC LIMIT 100" line.
$b4 = $ws.Range("B4").Text
$oldTail = "    order by f.file_name"
$newTail = "     order By f.file_name ASC LIMIT 100"
if ($b4.EndsWith($oldTail)) {
    $b4 = $b4.Substring(0, $b4.Length - $oldTail.Length) + $newTail
}
$ws.Range("B4").Value = $b4

# --- Row heights grew because the wrapped text in column B now spans one
#     more line (engine does not auto-measure wrapped text heights, so set
#     them explicitly to match Excel's own recalculated values).
$ws.Rows(2).RowHeight = 331.2
$ws.Rows(3).RowHeight = 360

# --- Selection moved from B4 to C4.
$ws.Range("C4").Select()
